$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.437.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.48%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.09%  "

$ws.Range("E6").Value = "  -2.29%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.22"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.92%  "

$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0617"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0877"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.889.13"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.566"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.439.43"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.79"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0727"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.54%  "

$ws.Range("E24").Value = "  -0.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.19"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.24%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  -2.28%  "

$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.448.36"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.73%  "

$ws.Range("E35").Value = "  -4.09%  "

$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.913"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.12%  "

$ws.Range("E38").Value = "  -4.57%  "

$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.05"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.45%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("E44").Value = "  -2.31%  "

$ws.Range("E45").Value = "  -2.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.797.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.37%  "

$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.53%  "

$ws.Range("E49").Value = "  -2.35%  "

$ws.Range("E50").Value = "  -1.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.13%  "
